# Add the new weekly price data for "Femacal de La Calera" (Ajo / Garlic).
# 4 new rows are inserted immediately before the old row 331, pushing the
# existing rows 331-351 down to 335-355 (dimension grows from R351 to R355).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("331:334").Insert()

# Row 331
$ws.Range("A331").Value = 3
$ws.Range("B331").Value = 'Femacal de La Calera'
$ws.Range("C331").Value = 'Coquimbo'
$ws.Range("D331").Value = 44585
$ws.Range("E331").Value = 5
$ws.Range("F331").Value = 100112003
$ws.Range("G331").Value = 'Ajo'
$ws.Range("H331").Value = 'Chino'
$ws.Range("I331").Value = '1a (cosecha)'
$ws.Range("J331").Value = 35
$ws.Range("K331").Value = 16000
$ws.Range("L331").Value = 16000
$ws.Range("M331").Value = 16000
$ws.Range("N331").Value = '$/caja 10 kilos'
$ws.Range("O331").Value = 'Llay Llay'
$ws.Range("P331").Value = 1600
$ws.Range("Q331").Value = 10
$ws.Range("R331").Value = 'Hortaliza'

# Row 332
$ws.Range("A332").Value = 3
$ws.Range("B332").Value = 'Femacal de La Calera'
$ws.Range("C332").Value = 'Coquimbo'
$ws.Range("D332").Value = 44585
$ws.Range("E332").Value = 5
$ws.Range("F332").Value = 100112003
$ws.Range("G332").Value = 'Ajo'
$ws.Range("H332").Value = 'Chino'
$ws.Range("I332").Value = '1a (cosecha)'
$ws.Range("J332").Value = 90
$ws.Range("K332").Value = 7000
$ws.Range("L332").Value = 7000
$ws.Range("M332").Value = 7000
$ws.Range("N332").Value = '$/trenza 50 unidades'
$ws.Range("O332").Value = 'Llay Llay'
$ws.Range("P332").Value = 1400
$ws.Range("Q332").Value = 5
$ws.Range("R332").Value = 'Hortaliza'

# Row 333
$ws.Range("A333").Value = 3
$ws.Range("B333").Value = 'Femacal de La Calera'
$ws.Range("C333").Value = 'Coquimbo'
$ws.Range("D333").Value = 44585
$ws.Range("E333").Value = 5
$ws.Range("F333").Value = 100112003
$ws.Range("G333").Value = 'Ajo'
$ws.Range("H333").Value = 'Chino'
$ws.Range("I333").Value = '2a (cosecha)'
$ws.Range("J333").Value = 80
$ws.Range("K333").Value = 5000
$ws.Range("L333").Value = 5000
$ws.Range("M333").Value = 5000
$ws.Range("N333").Value = '$/trenza 50 unidades'
$ws.Range("O333").Value = 'Llay Llay'
$ws.Range("P333").Value = 1000
$ws.Range("Q333").Value = 5
$ws.Range("R333").Value = 'Hortaliza'

# Row 334
$ws.Range("A334").Value = 3
$ws.Range("B334").Value = 'Femacal de La Calera'
$ws.Range("C334").Value = 'Coquimbo'
$ws.Range("D334").Value = 44585
$ws.Range("E334").Value = 5
$ws.Range("F334").Value = 100112003
$ws.Range("G334").Value = 'Ajo'
$ws.Range("H334").Value = 'Chino'
$ws.Range("I334").Value = 'Primera'
$ws.Range("J334").Value = 65
$ws.Range("K334").Value = 16000
$ws.Range("L334").Value = 16500
$ws.Range("M334").Value = 16231
$ws.Range("N334").Value = '$/caja 10 kilos'
$ws.Range("O334").Value = 'China'
$ws.Range("P334").Value = 1623
$ws.Range("Q334").Value = 10
$ws.Range("R334").Value = 'Hortaliza'
